# Apply the "add ground_warfae files and data" commit's data edits to the
# workbook. This updates cell values on the ship/SSM/SAM/inception sheets,
# strips the now-unused red/green highlight fills from the "ship" sheet,
# removes the custom column-E styling on "ship", sets a custom width for
# column A on "inception", and restores each sheet's stored selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "ship": updated figures + drop the red/green fill highlighting
# ---------------------------------------------------------------------
$shipWs = $wb.Worksheets.Item("ship")

$shipWs.Range("B7").Value  = 36
$shipWs.Range("C7").Value  = 36
$shipWs.Range("D7").Value  = 36
$shipWs.Range("E7").Value  = 36

$shipWs.Range("B10").Value = 15

$shipWs.Range("B12").Value = 15

$shipWs.Range("E13").Value = 1

$shipWs.Range("B14").Value = 1
$shipWs.Range("C14").Value = 5
$shipWs.Range("D14").Value = 5
$shipWs.Range("E14").Value = 5

$shipWs.Range("B24").Value = 150
$shipWs.Range("C24").Value = 150
$shipWs.Range("D24").Value = 150
$shipWs.Range("E24").Value = 150

$shipWs.Range("B34").Value = 2
$shipWs.Range("C34").Value = 2
$shipWs.Range("D34").Value = 2
$shipWs.Range("E34").Value = 2

$shipWs.Range("B35").Value = 3
$shipWs.Range("C35").Value = 3
$shipWs.Range("D35").Value = 3
$shipWs.Range("E35").Value = 3

# The workbook used to shade columns C/D (red) and E (green) via cell
# styles; that highlighting (and the dedicated column-E style) was removed.
$shipWs.Columns("E").ClearFormats()
$shipWs.Range("A1:E38").ClearFormats()

# ---------------------------------------------------------------------
# Sheet "SSM": updated figures
# ---------------------------------------------------------------------
$ssmWs = $wb.Worksheets.Item("SSM")

$ssmWs.Range("B2").Value = 2.5
$ssmWs.Range("C2").Value = 3
$ssmWs.Range("D2").Value = 4.5

$ssmWs.Range("B8").Value = 150
$ssmWs.Range("C8").Value = 150
$ssmWs.Range("D8").Value = 150

# ---------------------------------------------------------------------
# Sheet "SAM": updated figures
# ---------------------------------------------------------------------
$samWs = $wb.Worksheets.Item("SAM")

$samWs.Range("B2").Value = 5.5
$samWs.Range("C2").Value = 5.5
$samWs.Range("D2").Value = 5.5
$samWs.Range("E2").Value = 5.5

$samWs.Range("B3").Value = 60
$samWs.Range("C3").Value = 60
$samWs.Range("D3").Value = 30
$samWs.Range("E3").Value = 30

# ---------------------------------------------------------------------
# Sheet "patrol_aircraft": no data changes, just restore its selection
# ---------------------------------------------------------------------
$patrolWs = $wb.Worksheets.Item("patrol_aircraft")

# ---------------------------------------------------------------------
# Sheet "inception": updated figure + widen column A
# ---------------------------------------------------------------------
$inceptionWs = $wb.Worksheets.Item("inception")

$inceptionWs.Range("B2").Value = 90
$inceptionWs.Columns("A").ColumnWidth = 15.2857142857142857

# ---------------------------------------------------------------------
# Restore each sheet's saved selection (stored per-sheet regardless of
# which sheet ends up active). Do the currently-active "ship" sheet last
# so it remains the active tab/selection on save.
# ---------------------------------------------------------------------
$ssmWs.Range("C3").Select()
$samWs.Range("N11").Select()
$patrolWs.Range("F6").Select()
$inceptionWs.Range("B2").Select()

$shipWs.Activate()
$shipWs.Range("F6").Select()
